$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that were re-pulled / recalculated
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -5
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -9
$ws.Range("F7").Value = -5
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -4
